# Apply the data updates described in the commit diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Black Scholes" - update dates and valuation inputs
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Black Scholes")
$ws1.Range("B2:B5").NumberFormat = "@"
$ws1.Range("B2").Value = "1/10/2024"
$ws1.Range("B3").Value = "7/10/2024"
$ws1.Range("B4").Value = "7/10/2025"
$ws1.Range("B5").Value = "1/10/2025"
$ws1.Range("B6").Value = 150
$ws1.Range("B7").Value = 100
$ws1.Range("B8").Value = 1
$ws1.Range("B9").Value = 0.0522
$ws1.Range("B10").Value = 0.2323
$ws1.Range("B11").Value = 55.34

# ---------------------------------------------------------------
# Sheet "Volatility" - replace ticker list with a smaller set and
# recompute header/average; drop the now-unused trailing rows.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Volatility")
$ws2.Range("B1").Value = "2023-07-10 to 2024-07-10"
$ws2.Range("A2").Value = "AAPL"
$ws2.Range("B2").Value = 22.05
$ws2.Range("A3").Value = "GOOG"
$ws2.Range("B3").Value = 27.48
$ws2.Range("A4").Value = "MSFT"
$ws2.Range("B4").Value = 20.17
$ws2.Range("A5").Value = "Average"
$ws2.Range("B5").Value = 23.23
$ws2.Range("A6:B24").EntireRow.Delete()

# ---------------------------------------------------------------
# Sheet "Risk Free Rate" - update as-of date and yield curve values
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Risk Free Rate")
$ws3.Range("B1").NumberFormat = "@"
$ws3.Range("B1").Value = "2024-07-10"
$ws3.Range("B2").Value = 5.22
$ws3.Range("B3").Value = 4.97
$ws3.Range("B4").Value = 4.73
$ws3.Range("B5").Value = 4.49
$ws3.Range("B6").Value = 4.24
$ws3.Range("B7").Value = 4.25
$ws3.Range("B8").Value = 4.26
$ws3.Range("B9").Value = 4.26
$ws3.Range("B10").Value = 4.27
$ws3.Range("B11").Value = 4.28
$ws3.Range("B12").Value = 4.29
$ws3.Range("B13").Value = 4.3
$ws3.Range("B14").Value = 4.31
$ws3.Range("B15").Value = 4.32
$ws3.Range("B16").Value = 4.33
$ws3.Range("B17").Value = 4.34
$ws3.Range("B18").Value = 4.35
$ws3.Range("B19").Value = 4.36
$ws3.Range("B20").Value = 4.37
$ws3.Range("B21").Value = 4.38
$ws3.Range("B22").Value = 4.38
$ws3.Range("B23").Value = 4.39
$ws3.Range("B24").Value = 4.4
$ws3.Range("B25").Value = 4.41
$ws3.Range("B26").Value = 4.42
$ws3.Range("B27").Value = 4.43
$ws3.Range("B28").Value = 4.44
$ws3.Range("B29").Value = 4.45
$ws3.Range("B30").Value = 4.46
$ws3.Range("B31").Value = 4.47
